# Add a new phishing-sample row (row 29) to Sheet1, matching the
# "added new samples (thanks joseph)" commit.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$newRow = 29

# id
$ws.Cells.Item($newRow, 1).Value = 28
# type -> "msg" (already an existing shared string, reused)
$ws.Cells.Item($newRow, 2).Value = "msg"

# added (date) -> copy the date-formatted style from the row above (s="1",
# numFmtId 14 "m/d/yyyy") so we don't create a duplicate custom number format,
# then set the actual date value.
$ws.Cells.Item($newRow - 1, 3).Copy()
$ws.Cells.Item($newRow, 3).PasteSpecial(-4122)
$ws.Cells.Item($newRow, 3).Value = "7/3/2021"

# source
$ws.Cells.Item($newRow, 4).Value = "MCAST"
# url
$ws.Cells.Item($newRow, 5).Value = "shortened"
# motivation
$ws.Cells.Item($newRow, 6).Value = "delivery"
# language
$ws.Cells.Item($newRow, 7).Value = "mt"
# personalised
$ws.Cells.Item($newRow, 8).Value = "no"
# description (brand new shared string)
$ws.Cells.Item($newRow, 9).Value = "click to accept a packet"

$excel.CutCopyMode = 0

# Match the author's final selection/view state from the diff.
$null = $ws.Range("G29").Select()
